$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Case_X123"
$ws.Range("B4").Value = "Case_X123"
$ws.Range("B5").Value = "Case_X123"
$ws.Range("B6").Value = "Case_X123"
$ws.Range("B7").Value = "Case_X123"
$ws.Range("B12").Value = "Case_X123"
$ws.Range("B13").Value = "Case_X123"

$ws.Range("C2").Value = "TTName"
$ws.Range("C11").Value = "TTName"

$ws.Range("C11").Select()
